$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.217.73"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.180.89"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "255.75"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +6.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.627"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "67.84"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.28%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.570"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +6.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.75"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.75"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0933"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.10"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +8.72%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.504.21"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.870"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.46"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.172.38"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.150.17"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0953"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.62%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.80"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.23"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.04"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.93"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +9.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.86"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +22.79%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.55"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +6.71%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.92"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.48%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.123"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.63%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0742"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +6.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.44"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +6.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "27.20"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +17.26%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +10.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.61"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0300"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +13.37%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.55"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +22.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.68"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.60"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.09"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.200"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.63"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.09%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.06%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.29"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.43%  "
